$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reset Group A rows (2-5) -- Matches, Won, Lost, Points all back to 0
$ws.Range("B2:E5").Value = 0

# Group B row 7 (Rajastan_Australia) -> played 1 match, won it: Points = 2
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 2

# Group B rows 8-9 (Kolkata_England, Punjab_Pakistan) reset to 0
$ws.Range("B8:E9").Value = 0

# Group B row 10 (Sunrisers_SriLanka) -> played 1 match, lost it: Points = 0
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Append trailing whitespace to the three shared-string team names (shared strings idx 11-13)
$ws.Range("A2").Value = "Mumbai_India                         "
$ws.Range("A3").Value = "Chennai_SouthAfrica              "
$ws.Range("A4").Value = "Delhi_NewZealand                  "

# Update the current selection to match the recorded view state
$ws.Range("B7:E9").Select()
